# "added selenium grid code"
#
# shoppingcartdata (sheet4): the stray "x" value in H7 is removed, which
# shrinks the sheet's used range back down to A1:G2 and the selection is
# moved to A3:XFD15.
#
# registrationdatawithemail (sheet6): the four automation3 mailbox cells
# (C1:C4) are renamed to automation4, and the C1 hyperlink's cached
# "display" text (which had drifted to the automation1 address) is
# cleared so it falls back to showing the cell's own text; the cell also
# picks up the builtin Hyperlink style that its sibling E2 already had.

$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("shoppingcartdata")
$ws6 = $wb.Worksheets.Item("registrationdatawithemail")

# --- shoppingcartdata: drop the lone H7 cell ---------------------------
$ws4.Range("H7").ClearContents()

# Update the sheet's own selection without stealing the workbook's active
# tab away from registrationdatawithemail (restored right after).
[void]$ws4.Range("A3:XFD15").Select()
[void]$ws6.Activate()

# --- registrationdatawithemail: automation3 -> automation4 emails ------
$ws6.Range("C1").Value = "harpreet.automation4@gmail.com"
$ws6.Range("C2").Value = "siya.automation4@gmail.com"
$ws6.Range("C3").Value = "kriya.automation4@gmail.com"
$ws6.Range("C4").Value = "siyaa.automation4@gmail.com"

# Clear the stale cached hyperlink display text on C1 so the <hyperlink>
# element no longer carries a display="..." attribute.
$hyperlinks = $ws6.Hyperlinks
foreach ($link in $hyperlinks) {
    if ($link.Range.Address() -eq '$C$1') {
        $link.TextToDisplay = ""
    }
}

# Clearing TextToDisplay can reset the cell text, so reassert the value,
# then give C1 the same built-in Hyperlink style used elsewhere (e.g. E2).
$ws6.Range("C1").Value = "harpreet.automation4@gmail.com"
$ws6.Range("C1").Style = "Hyperlink"
